$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to snake_case field names ---
$ws.Cells.Item(1, 1).Value = 'mx_state'
$ws.Cells.Item(1, 2).Value = 'mx_municipality'
$ws.Cells.Item(1, 3).Value = 'n_matriculas'
$ws.Cells.Item(1, 4).Value = 'pct_matriculas'

# --- Normalize Spanish connector words (de/del/la/las/los/el/y) to title case ---
# in the state (col A) and municipality (col B) name columns
$ws.Cells.Item(7, 2).Value = 'Pabellón De Arteaga'
$ws.Cells.Item(8, 2).Value = 'Rincón De Romos'
$ws.Cells.Item(9, 2).Value = 'San Francisco De Los Romo'
$ws.Cells.Item(37, 2).Value = 'Amatenango De La Frontera'
$ws.Cells.Item(38, 2).Value = 'Amatenango Del Valle'
$ws.Cells.Item(42, 2).Value = 'Bejucal De Ocampo'
$ws.Cells.Item(44, 2).Value = 'Benemérito De Las Américas'
$ws.Cells.Item(53, 2).Value = 'Chiapa De Corzo'
$ws.Cells.Item(60, 2).Value = 'Comitán De Domínguez'
$ws.Cells.Item(91, 2).Value = 'Marqués De Comillas'
$ws.Cells.Item(92, 2).Value = 'Mazapa De Madero'
$ws.Cells.Item(97, 2).Value = 'Montecristo De Guerrero'
$ws.Cells.Item(101, 2).Value = 'Ocozocoautla De Espinosa'
$ws.Cells.Item(113, 2).Value = 'Salto De Agua'
$ws.Cells.Item(115, 2).Value = 'San Cristóbal De Las Casas'
$ws.Cells.Item(154, 2).Value = 'Guadalupe Y Calvo'
$ws.Cells.Item(157, 2).Value = 'Hidalgo Del Parral'
$ws.Cells.Item(164, 2).Value = 'San Francisco De Borja'
$ws.Cells.Item(165, 2).Value = 'San Francisco Del Oro'
$ws.Cells.Item(170, 1).Value = 'Ciudad De México'
$ws.Cells.Item(173, 2).Value = 'Cuajimalpa De Morelos'
$ws.Cells.Item(188, 1).Value = 'Coahuila De Zaragoza'
$ws.Cells.Item(204, 2).Value = 'San Juan De Sabinas'
$ws.Cells.Item(233, 2).Value = 'Nombre De Dios'
$ws.Cells.Item(239, 2).Value = 'Pánuco De Coronado'
$ws.Cells.Item(242, 2).Value = 'San Juan Del Río'
$ws.Cells.Item(249, 1).Value = 'Estado De México'
$ws.Cells.Item(249, 2).Value = 'Acambay De Ruíz Castañeda'
$ws.Cells.Item(252, 2).Value = 'Almoloya De Alquisiras'
$ws.Cells.Item(253, 2).Value = 'Almoloya De Juárez'
$ws.Cells.Item(259, 2).Value = 'Atizapán De Zaragoza'
$ws.Cells.Item(265, 2).Value = 'Chapa De Mota'
$ws.Cells.Item(267, 2).Value = 'Coacalco De Berriozábal'
$ws.Cells.Item(274, 2).Value = 'Ecatepec De Morelos'
$ws.Cells.Item(281, 2).Value = 'Ixtapan De La Sal'
$ws.Cells.Item(295, 2).Value = 'Naucalpan De Juárez'
$ws.Cells.Item(305, 2).Value = 'San Felipe Del Progreso'
$ws.Cells.Item(306, 2).Value = 'San Martín De Las Pirámides'
$ws.Cells.Item(308, 2).Value = 'San Simón De Guerrero'
$ws.Cells.Item(310, 2).Value = 'Soyaniquilpan De Juárez'
$ws.Cells.Item(319, 2).Value = 'Tenango Del Valle'
$ws.Cells.Item(329, 2).Value = 'Tlalnepantla De Baz'
$ws.Cells.Item(335, 2).Value = 'Valle De Bravo'
$ws.Cells.Item(336, 2).Value = 'Valle De Chalco Solidaridad'
$ws.Cells.Item(339, 2).Value = 'Villa De Allende'
$ws.Cells.Item(340, 2).Value = 'Villa Del Carbón'
$ws.Cells.Item(350, 2).Value = 'Apaseo El Alto'
$ws.Cells.Item(351, 2).Value = 'Apaseo El Grande'
$ws.Cells.Item(360, 2).Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Cells.Item(364, 2).Value = 'Jaral Del Progreso'
$ws.Cells.Item(371, 2).Value = 'Purísima Del Rincón'
$ws.Cells.Item(376, 2).Value = 'San Diego De La Unión'
$ws.Cells.Item(378, 2).Value = 'San Francisco Del Rincón'
$ws.Cells.Item(380, 2).Value = 'San Luis De La Paz'
$ws.Cells.Item(381, 2).Value = 'San Miguel De Allende'
$ws.Cells.Item(383, 2).Value = 'Santa Cruz De Juventino Rosas'
$ws.Cells.Item(384, 2).Value = 'Silao De La Victoria'
$ws.Cells.Item(389, 2).Value = 'Valle De Santiago'
$ws.Cells.Item(395, 2).Value = 'Acapulco De Juárez'
$ws.Cells.Item(398, 2).Value = 'Ajuchitlán Del Progreso'
$ws.Cells.Item(399, 2).Value = 'Alcozauca De Guerrero'
$ws.Cells.Item(403, 2).Value = 'Atenango Del Río'
$ws.Cells.Item(405, 2).Value = 'Atoyac De Álvarez'
$ws.Cells.Item(406, 2).Value = 'Ayutla De Los Libres'
$ws.Cells.Item(408, 2).Value = 'Chilapa De Álvarez'
$ws.Cells.Item(409, 2).Value = 'Chilpancingo De Los Bravo'
$ws.Cells.Item(410, 2).Value = 'Cochoapa El Grande'
$ws.Cells.Item(415, 2).Value = 'Coyuca De Benítez'
$ws.Cells.Item(416, 2).Value = 'Coyuca De Catalán'
$ws.Cells.Item(420, 2).Value = 'Cuetzala Del Progreso'
$ws.Cells.Item(421, 2).Value = 'Cutzamala De Pinzón'
$ws.Cells.Item(428, 2).Value = 'Huitzuco De Los Figueroa'
$ws.Cells.Item(429, 2).Value = 'Iguala De La Independencia'
$ws.Cells.Item(431, 2).Value = 'Ixcateopan De Cuauhtémoc'
$ws.Cells.Item(432, 2).Value = 'José Joaquín De Herrera'
$ws.Cells.Item(435, 2).Value = 'La Unión De Isidoro Montes De Oca'
$ws.Cells.Item(441, 2).Value = 'Mártir De Cuilapan'
$ws.Cells.Item(452, 2).Value = 'Taxco De Alarcón'
$ws.Cells.Item(455, 2).Value = 'Tepecoacuilco De Trujano'
$ws.Cells.Item(457, 2).Value = 'Tixtla De Guerrero'
$ws.Cells.Item(461, 2).Value = 'Tlalixtaquilla De Maldonado'
$ws.Cells.Item(462, 2).Value = 'Tlapa De Comonfort'
$ws.Cells.Item(464, 2).Value = 'Técpan De Galeana'
$ws.Cells.Item(469, 2).Value = 'Zihuatanejo De Azueta'
$ws.Cells.Item(480, 2).Value = 'Atotonilco El Grande'
$ws.Cells.Item(486, 2).Value = 'Cuautepec De Hinojosa'
$ws.Cells.Item(491, 2).Value = 'Huasca De Ocampo'
$ws.Cells.Item(495, 2).Value = 'Huejutla De Reyes'
$ws.Cells.Item(499, 2).Value = 'Jacala De Ledezma'
$ws.Cells.Item(506, 2).Value = 'Mineral Del Chico'
$ws.Cells.Item(507, 2).Value = 'Mineral Del Monte'
$ws.Cells.Item(508, 2).Value = 'Mixquiahuala De Juárez'
$ws.Cells.Item(509, 2).Value = 'Molango De Escamilla'
$ws.Cells.Item(511, 2).Value = 'Nopala De Villagrán'
$ws.Cells.Item(512, 2).Value = 'Omitlán De Juárez'
$ws.Cells.Item(513, 2).Value = 'Pachuca De Soto'
$ws.Cells.Item(516, 2).Value = 'Progreso De Obregón'
$ws.Cells.Item(522, 2).Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Cells.Item(523, 2).Value = 'Santiago De Anaya'
$ws.Cells.Item(527, 2).Value = 'Tenango De Doria'
$ws.Cells.Item(529, 2).Value = 'Tepehuacán De Guerrero'
$ws.Cells.Item(530, 2).Value = 'Tepeji Del Río De Ocampo'
$ws.Cells.Item(532, 2).Value = 'Tezontepec De Aldama'
$ws.Cells.Item(539, 2).Value = 'Tula De Allende'
$ws.Cells.Item(540, 2).Value = 'Tulancingo De Bravo'
$ws.Cells.Item(543, 2).Value = 'Zacualtipán De Ángeles'
$ws.Cells.Item(544, 2).Value = 'Zapotlán De Juárez'
$ws.Cells.Item(549, 2).Value = 'Acatlán De Juárez'
$ws.Cells.Item(550, 2).Value = 'Ahualulco De Mercado'
$ws.Cells.Item(555, 2).Value = 'Atotonilco El Alto'
$ws.Cells.Item(556, 2).Value = 'Autlán De Navarro'
$ws.Cells.Item(565, 2).Value = 'Cuautitlán De García Barragán'
$ws.Cells.Item(571, 2).Value = 'Encarnación De Díaz'
$ws.Cells.Item(575, 2).Value = 'Huejuquilla El Alto'
$ws.Cells.Item(578, 2).Value = 'Jilotlán De Los Dolores'
$ws.Cells.Item(582, 2).Value = 'Lagos De Moreno'
$ws.Cells.Item(588, 2).Value = 'Ojuelos De Jalisco'
$ws.Cells.Item(593, 2).Value = 'San Diego De Alejandría'
$ws.Cells.Item(596, 2).Value = 'San Miguel El Alto'
$ws.Cells.Item(598, 2).Value = 'Santa María Del Oro'
$ws.Cells.Item(601, 2).Value = 'Tamazula De Gordiano'
$ws.Cells.Item(603, 2).Value = 'Techaluta De Montenegro'
$ws.Cells.Item(605, 2).Value = 'Teocuitatlán De Corona'
$ws.Cells.Item(606, 2).Value = 'Tepatitlán De Morelos'
$ws.Cells.Item(608, 2).Value = 'Tlajomulco De Zúñiga'
$ws.Cells.Item(615, 2).Value = 'Unión De San Antonio'
$ws.Cells.Item(616, 2).Value = 'Unión De Tula'
$ws.Cells.Item(617, 2).Value = 'Valle De Juárez'
$ws.Cells.Item(620, 2).Value = 'Yahualica De González Gallo'
$ws.Cells.Item(621, 2).Value = 'Zacoalco De Torres'
$ws.Cells.Item(624, 2).Value = 'Zapotitlán De Vadillo'
$ws.Cells.Item(626, 2).Value = 'Zapotlán El Grande'
$ws.Cells.Item(628, 1).Value = 'Michoacán De Ocampo'
$ws.Cells.Item(646, 2).Value = 'Coalcomán De Vázquez Pallares'
$ws.Cells.Item(708, 2).Value = 'Tiquicheo De Nicolás Romero'
$ws.Cells.Item(728, 2).Value = 'Coatlán Del Río'
$ws.Cells.Item(740, 2).Value = 'Puente De Ixtla'
$ws.Cells.Item(745, 2).Value = 'Tetela Del Volcán'
$ws.Cells.Item(747, 2).Value = 'Tlaltizapán De Zapata'
$ws.Cells.Item(756, 2).Value = 'Amatlán De Cañas'
$ws.Cells.Item(757, 2).Value = 'Bahía De Banderas'
$ws.Cells.Item(760, 2).Value = 'Ixtlán Del Río'
$ws.Cells.Item(765, 2).Value = 'Santa María Del Oro'
$ws.Cells.Item(783, 2).Value = 'Mier Y Noriega'
$ws.Cells.Item(787, 2).Value = 'San Nicolás De Los Garza'
$ws.Cells.Item(792, 2).Value = 'Acatlán De Pérez Figueroa'
$ws.Cells.Item(800, 2).Value = 'Chalcatongo De Hidalgo'
$ws.Cells.Item(802, 2).Value = 'Coicoyán De Las Flores'
$ws.Cells.Item(805, 2).Value = 'Constancia Del Rosario'
$ws.Cells.Item(807, 2).Value = 'Cuilápam De Guerrero'
$ws.Cells.Item(809, 2).Value = 'El Barrio De La Soledad'
$ws.Cells.Item(810, 2).Value = 'Guadalupe De Ramírez'
$ws.Cells.Item(811, 2).Value = 'Guevea De Humboldt'
$ws.Cells.Item(812, 2).Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Cells.Item(813, 2).Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Cells.Item(814, 2).Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Cells.Item(815, 2).Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Cells.Item(816, 2).Value = 'Huajuapan De León'
$ws.Cells.Item(817, 2).Value = 'Huautla De Jiménez'
$ws.Cells.Item(818, 2).Value = 'Ixtlán De Juárez'
$ws.Cells.Item(822, 2).Value = 'Mariscala De Juárez'
$ws.Cells.Item(824, 2).Value = 'Miahuatlán De Porfirio Díaz'
$ws.Cells.Item(826, 2).Value = 'Mártires De Tacubaya'
$ws.Cells.Item(827, 2).Value = 'Nejapa De Madero'
$ws.Cells.Item(829, 2).Value = 'Oaxaca De Juárez'
$ws.Cells.Item(830, 2).Value = 'Ocotlán De Morelos'
$ws.Cells.Item(831, 2).Value = 'Pinotepa De Don Luis'
$ws.Cells.Item(833, 2).Value = 'Putla Villa De Guerrero'
$ws.Cells.Item(844, 2).Value = 'San Antonino El Alto'
$ws.Cells.Item(857, 2).Value = 'San Dionisio Del Mar'
$ws.Cells.Item(859, 2).Value = 'San Felipe Jalapa De Díaz'
$ws.Cells.Item(865, 2).Value = 'San Francisco Del Mar'
$ws.Cells.Item(899, 2).Value = 'San Juan Del Río'
$ws.Cells.Item(914, 2).Value = 'San Mateo Del Mar'
$ws.Cells.Item(926, 2).Value = 'San Miguel Del Puerto'
$ws.Cells.Item(927, 2).Value = 'San Miguel El Grande'
$ws.Cells.Item(946, 2).Value = 'San Pedro El Alto'
$ws.Cells.Item(947, 2).Value = 'San Pedro Y San Pablo Ayutla'
$ws.Cells.Item(948, 2).Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Cells.Item(949, 2).Value = 'San Pedro Y San Pablo Tequixtepec'
$ws.Cells.Item(969, 2).Value = 'Santa Inés Del Monte'
$ws.Cells.Item(980, 2).Value = 'Santa María Jalapa Del Marqués'
$ws.Cells.Item(993, 2).Value = 'Santa María Del Rosario'
$ws.Cells.Item(1026, 2).Value = 'Santo Domingo De Morelos'
$ws.Cells.Item(1034, 2).Value = 'Tataltepec De Valdés'
$ws.Cells.Item(1035, 2).Value = 'Teotitlán De Flores Magón'
$ws.Cells.Item(1037, 2).Value = 'Tezoatlán De Segura Y Luna'
$ws.Cells.Item(1038, 2).Value = 'Tlacolula De Matamoros'
$ws.Cells.Item(1039, 2).Value = 'Totontepec Villa De Morelos'
$ws.Cells.Item(1043, 2).Value = 'Villa Sola De Vega'
$ws.Cells.Item(1044, 2).Value = 'Villa De Chilapa De Díaz'
$ws.Cells.Item(1045, 2).Value = 'Villa De Etla'
$ws.Cells.Item(1046, 2).Value = 'Villa De Tututepec'
$ws.Cells.Item(1047, 2).Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Cells.Item(1048, 2).Value = 'Villa De Zaachila'
$ws.Cells.Item(1049, 2).Value = 'Yutanduchi De Guerrero'
$ws.Cells.Item(1052, 2).Value = 'Zimatlán De Álvarez'
$ws.Cells.Item(1068, 2).Value = 'Ayotoxco De Guerrero'
$ws.Cells.Item(1070, 2).Value = 'Chalchicomula De Sesma'
$ws.Cells.Item(1085, 2).Value = 'Cuetzalan Del Progreso'
$ws.Cells.Item(1098, 2).Value = 'Huehuetlán El Chico'
$ws.Cells.Item(1103, 2).Value = 'Huitzilan De Serdán'
$ws.Cells.Item(1104, 2).Value = 'Ixcamilpa De Guerrero'
$ws.Cells.Item(1107, 2).Value = 'Izúcar De Matamoros'
$ws.Cells.Item(1116, 2).Value = 'Los Reyes De Juárez'
$ws.Cells.Item(1121, 2).Value = 'Palmar De Bravo'
$ws.Cells.Item(1137, 2).Value = 'San Nicolás De Los Ranchos'
$ws.Cells.Item(1140, 2).Value = 'San Salvador El Seco'
$ws.Cells.Item(1141, 2).Value = 'San Salvador El Verde'
$ws.Cells.Item(1145, 2).Value = 'Tecali De Herrera'
$ws.Cells.Item(1151, 2).Value = 'Tepatlaxco De Hidalgo'
$ws.Cells.Item(1155, 2).Value = 'Tepexi De Rodríguez'
$ws.Cells.Item(1157, 2).Value = 'Tetela De Ocampo'
$ws.Cells.Item(1158, 2).Value = 'Teteles De Avila Castillo'
$ws.Cells.Item(1162, 2).Value = 'Tlacotepec De Benito Juárez'
$ws.Cells.Item(1174, 2).Value = 'Xayacatlán De Bravo'
$ws.Cells.Item(1190, 2).Value = 'Amealco De Bonfil'
$ws.Cells.Item(1192, 2).Value = 'Cadereyta De Montes'
$ws.Cells.Item(1199, 2).Value = 'Jalpan De Serra'
$ws.Cells.Item(1200, 2).Value = 'Landa De Matamoros'
$ws.Cells.Item(1203, 2).Value = 'Pinal De Amoles'
$ws.Cells.Item(1206, 2).Value = 'San Juan Del Río'
$ws.Cells.Item(1219, 2).Value = 'Armadillo De Los Infante'
$ws.Cells.Item(1220, 2).Value = 'Axtla De Terrazas'
$ws.Cells.Item(1223, 2).Value = 'Cerro De San Pedro'
$ws.Cells.Item(1227, 2).Value = 'Ciudad Del Maíz'
$ws.Cells.Item(1237, 2).Value = 'Mexquitic De Carmona'
$ws.Cells.Item(1242, 2).Value = 'San Ciro De Acosta'
$ws.Cells.Item(1248, 2).Value = 'Santa María Del Río'
$ws.Cells.Item(1249, 2).Value = 'Soledad De Graciano Sánchez'
$ws.Cells.Item(1255, 2).Value = 'Tanquián De Escobedo'
$ws.Cells.Item(1257, 2).Value = 'Villa De Arista'
$ws.Cells.Item(1258, 2).Value = 'Villa De Arriaga'
$ws.Cells.Item(1259, 2).Value = 'Villa De Guadalupe'
$ws.Cells.Item(1260, 2).Value = 'Villa De Ramos'
$ws.Cells.Item(1261, 2).Value = 'Villa De Reyes'
$ws.Cells.Item(1262, 2).Value = 'Villa De La Paz'
$ws.Cells.Item(1302, 2).Value = 'Jalpa De Méndez'
$ws.Cells.Item(1337, 2).Value = 'Soto La Marina'
$ws.Cells.Item(1348, 2).Value = 'Contla De Juan Cuamatzi'
$ws.Cells.Item(1354, 2).Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Cells.Item(1357, 2).Value = 'Papalotla De Xicohténcatl'
$ws.Cells.Item(1360, 2).Value = 'San Pablo Del Monte'
$ws.Cells.Item(1365, 2).Value = 'Tetla De La Solidaridad'
$ws.Cells.Item(1376, 1).Value = 'Veracruz De Ignacio De La Llave'
$ws.Cells.Item(1382, 2).Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Cells.Item(1385, 2).Value = 'Amatlán De Los Reyes'
$ws.Cells.Item(1395, 2).Value = 'Boca Del Río'
$ws.Cells.Item(1399, 2).Value = 'Castillo De Teayo'
$ws.Cells.Item(1401, 2).Value = 'Cazones De Herrera'
$ws.Cells.Item(1407, 2).Value = 'Chinampa De Gorostiza'
$ws.Cells.Item(1419, 2).Value = 'Cosamaloapan De Carpio'
$ws.Cells.Item(1437, 2).Value = 'Hueyapan De Ocampo'
$ws.Cells.Item(1438, 2).Value = 'Ignacio De La Llave'
$ws.Cells.Item(1442, 2).Value = 'Ixhuatlán De Madero'
$ws.Cells.Item(1443, 2).Value = 'Ixhuatlán Del Café'
$ws.Cells.Item(1444, 2).Value = 'Ixhuatlán Del Sureste'
$ws.Cells.Item(1452, 2).Value = 'Juchique De Ferrer'
$ws.Cells.Item(1457, 2).Value = 'Lerdo De Tejada'
$ws.Cells.Item(1462, 2).Value = 'Martínez De La Torre'
$ws.Cells.Item(1467, 2).Value = 'Mixtla De Altamirano'
$ws.Cells.Item(1469, 2).Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Cells.Item(1478, 2).Value = 'Ozuluama De Mascareñas'
$ws.Cells.Item(1480, 2).Value = 'Paso De Ovejas'
$ws.Cells.Item(1481, 2).Value = 'Paso Del Macho'
$ws.Cells.Item(1485, 2).Value = 'Poza Rica De Hidalgo'
$ws.Cells.Item(1494, 2).Value = 'Sayula De Alemán'
$ws.Cells.Item(1495, 2).Value = 'Soledad De Doblado'
$ws.Cells.Item(1499, 2).Value = 'Tatahuicapan De Juárez'
$ws.Cells.Item(1527, 2).Value = 'Vega De Alatorre'
$ws.Cells.Item(1537, 2).Value = 'Zontecomatlán De López Y Fuentes'
$ws.Cells.Item(1538, 2).Value = 'Zozocolco De Hidalgo'
$ws.Cells.Item(1555, 2).Value = 'Cañitas De Felipe Pescador'
$ws.Cells.Item(1556, 2).Value = 'Concepción Del Oro'
$ws.Cells.Item(1565, 2).Value = 'Jiménez Del Teul'
$ws.Cells.Item(1573, 2).Value = 'Nochistlán De Mejía'
$ws.Cells.Item(1574, 2).Value = 'Noria De Ángeles'
$ws.Cells.Item(1582, 2).Value = 'Teúl De González Ortega'
$ws.Cells.Item(1583, 2).Value = 'Tlaltenango De Sánchez Román'

# --- Row 1591: "TOTAL" -> "Total" ---
$ws.Cells.Item(1591, 1).Value = 'Total'

# --- Remove trailing metadata/footer rows (1593:1597); dimension becomes A1:D1591 ---
$ws.Rows("1593:1597").Delete()

